# Apply "results with fixed workflow" update
# Sheet 1: iModulon — recompute BAR_Set_percentage / NBR_Set_percentage for each row (row count unchanged: 13 data rows)
# Sheet 2: Subsystem — full rebuild: new category rows inserted, all rows re-sorted/re-derived (13->? actually 52->56 data rows)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("iModulon")
$ws2 = $wb.Worksheets.Item("Subsystem")

$sheet1Data = @(
    @("Carbon Metabolism", 1.12781954887218, 0.0),
    @("Element_homeostasis", 1.691729323308271, 0.0),
    @("Extracellular_structures", 8.082706766917292, 10.81081081081081),
    @("Global_stress", 7.330827067669173, 10.81081081081081),
    @("Iron_starvation", 0.1879699248120301, 0.0),
    @("Light_response", 2.067669172932331, 0.0),
    @("Metal_homeostasis", 1.879699248120301, 0.0),
    @("Nitrogen_metabolism", 0.9398496240601504, 1.351351351351351),
    @("Photosynthetic_response", 15.97744360902256, 6.756756756756757),
    @("Recombination", 0.7518796992481203, 0.0),
    @("functional", 0.5639097744360901, 0.0),
    @("translation", 0.9398496240601504, 1.351351351351351),
    @("uncharacterized", 3.759398496240602, 4.054054054054054),
)

$sheet2Data = @(
    @("Alanine, aspartate and glutamate metabolism", 1.538461538461539, 1.694915254237288),
    @("Aminosugars metabolism", 0.8547008547008548, 0.847457627118644),
    @("Arginine and proline metabolism", 2.905982905982906, 10.16949152542373),
    @("Biomass and maintenance functions", 0.1709401709401709, 0.0),
    @("Biotin metabolism", 0.0, 2.542372881355933),
    @("C5-Branched dibasic acid metabolism", 0.5128205128205128, 0.0),
    @("Carbon fixation", 1.025641025641026, 0.0),
    @("Carotenoid Biosynthesis", 2.222222222222222, 0.0),
    @("Citrate cycle (TCA cycle)", 1.025641025641026, 0.0),
    @("Cyanophycin metabolism", 0.0, 0.0),
    @("Extracellular exchange", 3.247863247863248, 0.0),
    @("Fatty acid biosynthesis", 17.09401709401709, 0.0),
    @("Folate biosynthesis", 2.222222222222222, 1.694915254237288),
    @("Fructose and mannose metabolism", 0.0, 6.779661016949152),
    @("Galactolipids metabolism", 4.615384615384616, 0.0),
    @("Glutamate metabolism", 1.196581196581197, 0.0),
    @("Glutathione metabolism", 0.0, 1.694915254237288),
    @("Glycerolipid metabolism", 0.5128205128205128, 0.847457627118644),
    @("Glycolysis/Gluconeogenesis", 3.247863247863248, 4.23728813559322),
    @("Glyoxylate and dicarboxylate metabolism", 2.051282051282051, 0.847457627118644),
    @("Histidine metabolism", 1.538461538461539, 1.694915254237288),
    @("Hydrogen production", 0.0, 0.0),
    @("Inositol phosphate metabolism", 0.0, 0.847457627118644),
    @("Intracellular demand", 0.1709401709401709, 0.0),
    @("Intracellular source/sink", 0.1709401709401709, 0.0),
    @("Lipopolysaccharide biosynthesis", 0.8547008547008548, 0.0),
    @("Lysine metabolism", 1.538461538461539, 1.694915254237288),
    @("Nicotinate and nicotinamide metabolism", 0.8547008547008548, 1.694915254237288),
    @("Nitrogen metabolism", 1.88034188034188, 0.847457627118644),
    @("Nucleotide sugars metabolism", 0.0, 3.389830508474576),
    @("Others", 0.6837606837606838, 1.694915254237288),
    @("Oxidative phosphorylation", 1.025641025641026, 0.847457627118644),
    @("PHB byosynthesis", 0.0, 1.694915254237288),
    @("Pantothenate and CoA biosynthesis", 1.538461538461539, 0.0),
    @("Pentose phosphate pathway", 1.196581196581197, 0.847457627118644),
    @("Peptidoglycan biosynthesis", 1.367521367521368, 0.847457627118644),
    @("Phenylalanine tyrosine and tryptophan biosynthesis", 3.247863247863248, 5.932203389830509),
    @("Photosynthesis", 0.8547008547008548, 0.0),
    @("Porphyrin and chlorophyll metabolism", 7.008547008547009, 5.084745762711865),
    @("Purine metabolism", 4.273504273504273, 7.627118644067797),
    @("Pyrimidine metabolism", 3.931623931623931, 1.694915254237288),
    @("Pyruvate metabolism", 1.538461538461539, 0.847457627118644),
    @("Riboflavin metabolism", 1.70940170940171, 0.0),
    @("Starch and sucrose metabolism", 0.8547008547008548, 3.389830508474576),
    @("Steroid biosynthesis", 0.0, 2.542372881355933),
    @("Sterol biosynthesis", 0.8547008547008548, 0.0),
    @("Sulfolipid Biosynthesis", 1.70940170940171, 0.0),
    @("Sulfur Cysteine and methionine metabolism", 2.393162393162394, 4.23728813559322),
    @("Terpenoid backbone biosynthesis", 1.538461538461539, 0.0),
    @("Thiamine metabolism", 0.0, 4.23728813559322),
    @("Transport", 7.179487179487179, 8.47457627118644),
    @("Ubiquinone and other pterpenoids biosynthesis", 2.393162393162394, 0.0),
    @("Urea cycle and metabolism of amino groups", 1.196581196581197, 0.847457627118644),
    @("Valine leucine and isoleucine biosynthesis", 2.051282051282051, 3.389830508474576),
    @("Vitamin B6 metabolism", 0.0, 4.23728813559322),
)

# --- Sheet1 (iModulon): overwrite existing data rows 2..14 (A unchanged, B/C updated) ---
for ($i = 0; $i -lt $sheet1Data.Length; $i++) {
    $row = $sheet1Data[$i]
    $r = $i + 2
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
}

# --- Sheet2 (Subsystem): overwrite/extend data rows 2..56 ---
$oldLastRow = 52
for ($i = 0; $i -lt $sheet2Data.Length; $i++) {
    $row = $sheet2Data[$i]
    $r = $i + 2
    $isNewRow = $r -gt $oldLastRow

    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]

    if ($isNewRow) {
        # Copy column-A formatting (bold/border/center) from the last pre-existing row
        $ws2.Cells.Item($oldLastRow, 1).Copy()
        $ws2.Cells.Item($r, 1).PasteSpecial(-4122)
    }
}
$excel.CutCopyMode = $false
